$wb = $excel.ActiveWorkbook

# Add the new "simple" worksheet (a minimal key/value test fixture that
# exercises a leading "#Loc" style key without tripping up the kv
# transpiler anymore).
$ws = $wb.Worksheets.Add()
$ws.Name = "simple"

$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 2

$ws.Range("A2").Value = "key"
$ws.Range("B2").Value = "value"
$ws.Range("C2").Value = "#Loc_Test"

$ws.Range("A3").Value = "test"
$ws.Range("B3").Value = "test1"

$ws.Range("A4").Value = "test2"
$ws.Range("B4").Value = 1

# Put it right after "test_sheet1" (i.e. as the 2nd tab).
$ws.Move($null, $wb.Worksheets.Item("test_sheet1"))

# The move above can invalidate the old sheet handle (it's resolved by
# position), so re-resolve by name before touching the view state.
$simple = $wb.Worksheets.Item("simple")
$simple.Activate()
$simple.Range("D8").Select() | Out-Null
